$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.438.57'
$ws.Range('E2').Value = '  -1.14%  '

# Row 3
$ws.Range('D3').Value = '3.236.03'
$ws.Range('E3').Value = '  -1.24%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '579.04'
$ws.Range('E5').Value = '  -1.62%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '183.84'
$ws.Range('E6').Value = '  -1.21%  '

# Row 7
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.610'
$ws.Range('E8').Value = '  +1.43%  '

# Row 9
$ws.Range('D9').Value = '3.234.43'
$ws.Range('E9').Value = '  -1.24%  '

# Row 10
$ws.Range('E10').Value = '  -3.28%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.59'
$ws.Range('E11').Value = '  -2.19%  '

# Row 12
$ws.Range('E12').Value = '  -1.59%  '

# Row 13
$ws.Range('D13').Value = '3.797.18'
$ws.Range('E13').Value = '  -1.22%  '

# Row 14
$ws.Range('E14').Value = '  +0.04%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '27.62'
$ws.Range('E15').Value = '  -3.60%  '

# Row 16
$ws.Range('D16').Value = '67.496.29'
$ws.Range('E16').Value = '  -1.04%  '

# Row 17
$ws.Range('E17').Value = '  -1.98%  '

# Row 18
$ws.Range('D18').Value = '3.260.49'
$ws.Range('E18').Value = '  -0.32%  '

# Row 19
$ws.Range('E19').Value = '  -1.79%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.47'
$ws.Range('E20').Value = '  -1.38%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '395.15'
$ws.Range('E21').Value = '  +3.53%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.56'
$ws.Range('E22').Value = '  -2.41%  '

# Row 23
$ws.Range('E23').Value = '  -0.04%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '71.23'
$ws.Range('E24').Value = '  -0.39%  '

# Row 25
$ws.Range('E25').Value = '  +0.18%  '

# Row 26
$ws.Range('E26').Value = '  -2.75%  '

# Row 27
$ws.Range('E27').Value = '  -1.86%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.59'
$ws.Range('E28').Value = '  -2.00%  '

# Row 29
$ws.Range('E29').Value = '  +0.00%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.95'
$ws.Range('E30').Value = '  -2.19%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.54'
$ws.Range('E31').Value = '  -4.87%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '22.59'
$ws.Range('E32').Value = '  -1.61%  '

# Row 33
$ws.Range('E33').Value = '  -2.09%  '

# Row 34
$ws.Range('E34').Value = '  -2.58%  '

# Row 35
$ws.Range('E35').Value = '  +0.01%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '160.92'
$ws.Range('E36').Value = '  -1.42%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.47'
$ws.Range('E37').Value = '  -4.39%  '

# Row 38
$ws.Range('E38').Value = '  +0.79%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '26.47'
$ws.Range('E39').Value = '  -0.86%  '

# Row 40
$ws.Range('E40').Value = '  -4.52%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.55'
$ws.Range('E41').Value = '  -1.42%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.46'
$ws.Range('E42').Value = '  -5.09%  '

# Row 43
$ws.Range('E43').Value = '  -6.85%  '

# Row 44
$ws.Range('E44').Value = '  -0.94%  '

# Row 45
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '40.54'
$ws.Range('E45').Value = '  -2.14%  '

# Row 46
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.613.01'
$ws.Range('E46').Value = '  -1.12%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '24.66'
$ws.Range('E47').Value = '  -3.57%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '333.61'
$ws.Range('E48').Value = '  -2.94%  '

# Row 49
$ws.Range('E49').Value = '  -2.55%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.30'
$ws.Range('E50').Value = '  +0.83%  '

# Row 51
$ws.Range('E51').Value = '  -0.61%  '
